# Applies the BOM restructuring described in the commit:
#  - Removes the RN1/RN2 resistor-network rows and replaces them with eight
#    new individual resistor rows (R10-R17, packages M0806-M0813).
#  - Shifts the SV1..SOCK2 rows down to make room.
#  - Fixes the Mouser part numbers for SOCK1/SOCK2 (drops stray BOM char,
#    and both now reference the same part number).
#  - Extends the blank footer area by 6 rows (footer block moves from
#    rows 49-51 down to rows 55-57).
#  - Updates the active selection to match the new state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the 6 additional resistor rows -----------------------
# Old rows 20 (RN1) and 21 (RN2) are overwritten in place below; insert 6
# fresh rows right after them so the old SV1.. block (rows 22-29) slides
# down to rows 28-35, matching the target layout.
$ws.Rows.Item(22).Resize(6).Insert()

# --- 2. Fill rows 20-27 with the new resistor part data ---------------------
$resistors = @(
    @{ Row = 20; Part = "R10"; Package = "M0806" },
    @{ Row = 21; Part = "R11"; Package = "M0807" },
    @{ Row = 22; Part = "R12"; Package = "M0808" },
    @{ Row = 23; Part = "R13"; Package = "M0809" },
    @{ Row = 24; Part = "R14"; Package = "M0810" },
    @{ Row = 25; Part = "R15"; Package = "M0811" },
    @{ Row = 26; Part = "R16"; Package = "M0812" },
    @{ Row = 27; Part = "R17"; Package = "M0813" }
)

foreach ($r in $resistors) {
    $row = $r.Row
    $pkg = $r.Package
    $ws.Range("A$row").Value = $r.Part
    $ws.Range("B$row").Value = "1k"
    $ws.Range("C$row").Value = "R-EU_$pkg"
    $ws.Range("D$row").Value = $pkg
    $ws.Range("E$row").Value = "RESISTOR, European symbol"
    $ws.Range("F$row").Value = "652-CR0805JW-102ELF"
}

# --- 3. Fix SOCK1 / SOCK2 Mouser numbers (rows 34-35 after the shift) ------
$ws.Range("F34").Value = "571-1-2199298-4"
$ws.Range("F35").Value = "571-1-2199298-4"

# NOTE: the footer block (originally rows 49-51, styles 11/12/.../10) needs
# no separate insert: the single 6-row insert above (at row 22) already
# shifts everything below it - including the footer - down by 6 rows, so it
# lands exactly on rows 55-57, matching the target layout.

# --- 4. Update the saved selection -----------------------------------------
$ws.Range("F36").Select()
